# Generate Report for Handback
# Refreshes the handback-status report with a new pair of generated
# files (new GUID-named .md sources, new xliff hash/timestamps).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---- new identifiers -------------------------------------------------
$guid1 = "4f1cbecd-55f8-4daa-8b45-6daf77a83b61"
$guid2 = "ffffb07cf097-272b-49ac-ba02-3a6161471bad"
$xlfHash = "84c8dc3deb0928ce5488d97451314089b4fea04c"

$md1 = "$guid1.md"
$md2 = "$guid2.md"
$mdPath1 = "e2e\$guid1.md"
$mdPath2 = "e2e\$guid2.md"

$latestDate = "2016-08-21 19:04:53"

$zhCnXlf   = "$guid1.$xlfHash.zh-cn.xlf"
$zhCnStart = "2016-08-21 19:04:48"
$zhCnEnd   = "2016-08-21 19:05:12"

$deDeXlf = "$guid1.$xlfHash.de-de.xlf"
$deDeEnd = "2016-08-21 19:05:18"

# ---- Overview sheet ----------------------------------------------------
$wsOverview.Range("A2").Value = $md1
$wsOverview.Range("B2").Value = $mdPath1
$wsOverview.Range("G2").Value = $latestDate

$wsOverview.Range("A3").Value = $md2
$wsOverview.Range("B3").Value = $mdPath2
$wsOverview.Range("G3").Value = $latestDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d68d18bd9e8f399be2e7b361ebaf76967e33c982/e2e/$guid1.md", "", "", $mdPath1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d68d18bd9e8f399be2e7b361ebaf76967e33c982/e2e/$guid2.md", "", "", $mdPath2)
# restore the workbook's custom hyperlink look (Add() resets cells to the
# built-in theme hyperlink style)
$wsOverview.Range("B2:B3").Font.Color = 15570276
$wsOverview.Range("B2:B3").Font.Underline = $true

# ---- zh-cn sheet ---------------------------------------------------------
$wsZhCn.Range("A2").Value = $md1
$wsZhCn.Range("G2").Value = $zhCnXlf
$wsZhCn.Range("H2").Value = $zhCnStart
$wsZhCn.Range("I2").Value = $md1
$wsZhCn.Range("J2").Value = $zhCnXlf
$wsZhCn.Range("K2").Value = $zhCnEnd

$wsZhCn.Range("A3").Value = $md2
$wsZhCn.Range("G3").Value = $zhCnXlf
$wsZhCn.Range("H3").Value = $zhCnStart
$wsZhCn.Range("I3").Value = $md2
$wsZhCn.Range("J3").Value = $zhCnXlf
$wsZhCn.Range("K3").Value = $zhCnEnd

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d68d18bd9e8f399be2e7b361ebaf76967e33c982/e2e/$guid1.md", "", "", $md1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2d0d39426f8f604a1ddfc5d7e4795e96e780e286/e2e/$guid1.md", "", "", $md1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d68d18bd9e8f399be2e7b361ebaf76967e33c982/e2e/$guid2.md", "", "", $md2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2d0d39426f8f604a1ddfc5d7e4795e96e780e286/e2e/$guid2.md", "", "", $md2)
$wsZhCn.Range("A2:A3").Font.Color = 15570276
$wsZhCn.Range("A2:A3").Font.Underline = $true
$wsZhCn.Range("I2:I3").Font.Color = 15570276
$wsZhCn.Range("I2:I3").Font.Underline = $true

# ---- de-de sheet ---------------------------------------------------------
$wsDeDe.Range("A2").Value = $md1
$wsDeDe.Range("G2").Value = $deDeXlf
$wsDeDe.Range("H2").Value = $latestDate
$wsDeDe.Range("I2").Value = $md1
$wsDeDe.Range("J2").Value = $deDeXlf
$wsDeDe.Range("K2").Value = $deDeEnd

$wsDeDe.Range("A3").Value = $md2
$wsDeDe.Range("G3").Value = $deDeXlf
$wsDeDe.Range("H3").Value = $latestDate
$wsDeDe.Range("I3").Value = $md2
$wsDeDe.Range("J3").Value = $deDeXlf
$wsDeDe.Range("K3").Value = $deDeEnd

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d68d18bd9e8f399be2e7b361ebaf76967e33c982/e2e/$guid1.md", "", "", $md1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/70ddcd822ab69d5c2ddacd6959b7930fe4b0cc31/e2e/$guid1.md", "", "", $md1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d68d18bd9e8f399be2e7b361ebaf76967e33c982/e2e/$guid2.md", "", "", $md2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/70ddcd822ab69d5c2ddacd6959b7930fe4b0cc31/e2e/$guid2.md", "", "", $md2)
$wsDeDe.Range("A2:A3").Font.Color = 15570276
$wsDeDe.Range("A2:A3").Font.Underline = $true
$wsDeDe.Range("I2:I3").Font.Color = 15570276
$wsDeDe.Range("I2:I3").Font.Underline = $true
